$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 4000
$ws.Cells.Item(64, 9).Value = 4000
$ws.Cells.Item(64, 11).Value = 4000
$ws.Cells.Item(64, 13).Value = -3752
$ws.Cells.Item(67, 8).Value = 4000
$ws.Cells.Item(67, 9).Value = 4000
$ws.Cells.Item(67, 11).Value = 4000
$ws.Cells.Item(67, 13).Value = -3142
$ws.Cells.Item(69, 8).Value = 3677.5
$ws.Cells.Item(69, 9).Value = 3646.3333
$ws.Cells.Item(69, 10).Value = 3684.6924
$ws.Cells.Item(69, 11).Value = 10938.9999
$ws.Cells.Item(69, 12).Value = 11054.0772
$ws.Cells.Item(69, 13).Value = -10064.9999
$ws.Cells.Item(69, 14).Value = -12802.0772
$ws.Cells.Item(72, 8).Value = 3677.5
$ws.Cells.Item(72, 9).Value = 3646.3333
$ws.Cells.Item(72, 10).Value = 3684.6924
$ws.Cells.Item(72, 11).Value = 32816.9997
$ws.Cells.Item(72, 12).Value = 33162.2316
$ws.Cells.Item(72, 13).Value = -28448.9997
$ws.Cells.Item(72, 14).Value = -41898.2316
$ws.Cells.Item(76, 8).Value = 3323.0386
$ws.Cells.Item(76, 9).Value = 3304.3044
$ws.Cells.Item(76, 10).Value = 3466.6667
$ws.Cells.Item(76, 11).Value = 3304.3044
$ws.Cells.Item(76, 12).Value = 3466.6667
$ws.Cells.Item(76, 13).Value = -2989.3044
$ws.Cells.Item(76, 14).Value = -4096.6667
$ws.Cells.Item(79, 8).Value = 3323.0386
$ws.Cells.Item(79, 9).Value = 3304.3044
$ws.Cells.Item(79, 10).Value = 3466.6667
$ws.Cells.Item(79, 11).Value = 3304.3044
$ws.Cells.Item(79, 12).Value = 3466.6667
$ws.Cells.Item(79, 13).Value = -2212.3044
$ws.Cells.Item(79, 14).Value = -5650.6667
$ws.Cells.Item(96, 8).Value = 1536.75
$ws.Cells.Item(96, 9).Value = 800
$ws.Cells.Item(96, 10).Value = 1782.3334
$ws.Cells.Item(96, 11).Value = 2400
$ws.Cells.Item(96, 12).Value = 5347.0002
$ws.Cells.Item(96, 13).Value = -1027
$ws.Cells.Item(96, 14).Value = -8093.0002
$ws.Cells.Item(137, 8).Value = 4213.8184
$ws.Cells.Item(137, 9).Value = 2745.2104
$ws.Cells.Item(137, 10).Value = 13515
$ws.Cells.Item(137, 11).Value = 8235.6312
$ws.Cells.Item(137, 12).Value = 40545
$ws.Cells.Item(137, 13).Value = -5685.6312
$ws.Cells.Item(137, 14).Value = -45645
$ws.Cells.Item(138, 8).Value = 2024.4
$ws.Cells.Item(138, 9).Value = 1945.2727
$ws.Cells.Item(138, 10).Value = 2100.087
$ws.Cells.Item(138, 11).Value = 5835.8181
$ws.Cells.Item(138, 12).Value = 6300.261
$ws.Cells.Item(138, 13).Value = -695.8181000000004
$ws.Cells.Item(138, 14).Value = -16580.261

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(23, 8).Value = 12675.2
$ws.Cells.Item(23, 10).Value = 9657.789000000001
$ws.Cells.Item(23, 12).Value = 9657.789000000001
$ws.Cells.Item(23, 14).Value = -10175.789
$ws.Cells.Item(32, 8).Value = 5439.43
$ws.Cells.Item(32, 9).Value = 3590.1667
$ws.Cells.Item(32, 10).Value = 22082.8
$ws.Cells.Item(32, 11).Value = 3590.1667
$ws.Cells.Item(32, 12).Value = 22082.8
$ws.Cells.Item(32, 13).Value = -3303.1667
$ws.Cells.Item(32, 14).Value = -22656.8
$ws.Cells.Item(37, 8).Value = 11943
$ws.Cells.Item(37, 10).Value = 11943
$ws.Cells.Item(37, 12).Value = 11943
$ws.Cells.Item(37, 14).Value = -12489
$ws.Cells.Item(44, 8).Value = 15498.333
$ws.Cells.Item(44, 10).Value = 15498.333
$ws.Cells.Item(44, 12).Value = 15498.333
$ws.Cells.Item(44, 14).Value = -16474.333
$ws.Cells.Item(55, 8).Value = 29999
$ws.Cells.Item(55, 10).Value = 29999
$ws.Cells.Item(55, 12).Value = 29999
$ws.Cells.Item(55, 14).Value = -30629
$ws.Cells.Item(63, 8).Value = 5777.4287
$ws.Cells.Item(63, 10).Value = 7712.2856
$ws.Cells.Item(63, 12).Value = 7712.2856
$ws.Cells.Item(63, 14).Value = -9084.285599999999
$ws.Cells.Item(66, 8).Value = 5777.4287
$ws.Cells.Item(66, 10).Value = 7712.2856
$ws.Cells.Item(66, 12).Value = 38561.428
$ws.Cells.Item(66, 14).Value = -45425.428
$ws.Cells.Item(75, 8).Value = 32000
$ws.Cells.Item(75, 10).Value = 32000
$ws.Cells.Item(75, 12).Value = 32000
$ws.Cells.Item(75, 14).Value = -33748
$ws.Cells.Item(78, 8).Value = 32000
$ws.Cells.Item(78, 10).Value = 32000
$ws.Cells.Item(78, 12).Value = 96000
$ws.Cells.Item(78, 14).Value = -104736
$ws.Cells.Item(80, 8).Value = 21427.572
$ws.Cells.Item(80, 10).Value = 21427.572
$ws.Cells.Item(80, 12).Value = 21427.572
$ws.Cells.Item(80, 14).Value = -23423.572
$ws.Cells.Item(83, 8).Value = 21427.572
$ws.Cells.Item(83, 10).Value = 21427.572
$ws.Cells.Item(83, 12).Value = 64282.716
$ws.Cells.Item(83, 14).Value = -74266.716
$ws.Cells.Item(110, 8).Value = 2122.5
$ws.Cells.Item(110, 9).Value = 2000
$ws.Cells.Item(110, 10).Value = 2245
$ws.Cells.Item(110, 11).Value = 2000
$ws.Cells.Item(110, 12).Value = 2245
$ws.Cells.Item(110, 13).Value = 45
$ws.Cells.Item(110, 14).Value = -6335

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 1131.4
$ws.Cells.Item(94, 9).Value = 798.4
$ws.Cells.Item(94, 10).Value = 1297.9
$ws.Cells.Item(94, 11).Value = 798.4
$ws.Cells.Item(94, 12).Value = 1297.9
$ws.Cells.Item(94, 13).Value = -347.4
$ws.Cells.Item(94, 14).Value = -2199.9
$ws.Cells.Item(134, 8).Value = 1961.931
$ws.Cells.Item(134, 9).Value = 1283.9524
$ws.Cells.Item(134, 10).Value = 3741.625
$ws.Cells.Item(134, 11).Value = 3851.857199999999
$ws.Cells.Item(134, 12).Value = 11224.875
$ws.Cells.Item(134, 13).Value = -1316.857199999999
$ws.Cells.Item(134, 14).Value = -16294.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(109, 8).Value = 3785.4546
$ws.Cells.Item(109, 9).Value = 1125
$ws.Cells.Item(109, 10).Value = 5305.7144
$ws.Cells.Item(109, 11).Value = 3375
$ws.Cells.Item(109, 12).Value = 15917.1432
$ws.Cells.Item(109, 13).Value = -2335
$ws.Cells.Item(109, 14).Value = -17997.1432
$ws.Cells.Item(131, 8).Value = 1080.7954
$ws.Cells.Item(131, 9).Value = 915
$ws.Cells.Item(131, 10).Value = 1088.6904
$ws.Cells.Item(131, 11).Value = 2745
$ws.Cells.Item(131, 12).Value = 3266.0712
$ws.Cells.Item(131, 13).Value = 2295
$ws.Cells.Item(131, 14).Value = -13346.0712
$ws.Cells.Item(137, 8).Value = 7584019.5
$ws.Cells.Item(137, 9).Value = 20848778
$ws.Cells.Item(137, 10).Value = 4157.9287
$ws.Cells.Item(137, 11).Value = 62546334
$ws.Cells.Item(137, 12).Value = 12473.7861
$ws.Cells.Item(137, 13).Value = -62541234
$ws.Cells.Item(137, 14).Value = -22673.7861

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 50942904
$ws.Cells.Item(80, 9).Value = 72717576
$ws.Cells.Item(80, 10).Value = 135333.33
$ws.Cells.Item(80, 11).Value = 72717576
$ws.Cells.Item(80, 12).Value = 135333.33
$ws.Cells.Item(80, 13).Value = -72716578
$ws.Cells.Item(80, 14).Value = -137329.33
$ws.Cells.Item(83, 8).Value = 50942904
$ws.Cells.Item(83, 9).Value = 72717576
$ws.Cells.Item(83, 10).Value = 135333.33
$ws.Cells.Item(83, 11).Value = 363587880
$ws.Cells.Item(83, 12).Value = 676666.6499999999
$ws.Cells.Item(83, 13).Value = -363582888
$ws.Cells.Item(83, 14).Value = -686650.6499999999
$ws.Cells.Item(132, 8).Value = 3275.5
$ws.Cells.Item(132, 9).Value = 3100.875
$ws.Cells.Item(132, 10).Value = 3624.75
$ws.Cells.Item(132, 11).Value = 9302.625
$ws.Cells.Item(132, 12).Value = 10874.25
$ws.Cells.Item(132, 13).Value = -6772.625
$ws.Cells.Item(132, 14).Value = -15934.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 1731.5
$ws.Cells.Item(68, 9).Value = 1731.5
$ws.Cells.Item(68, 11).Value = 1731.5
$ws.Cells.Item(68, 13).Value = -982.5
$ws.Cells.Item(71, 8).Value = 1731.5
$ws.Cells.Item(71, 9).Value = 1731.5
$ws.Cells.Item(71, 11).Value = 8657.5
$ws.Cells.Item(71, 13).Value = -4913.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 41217.5
$ws.Cells.Item(62, 9).Value = 2400
$ws.Cells.Item(62, 10).Value = 54156.668
$ws.Cells.Item(62, 11).Value = 2400
$ws.Cells.Item(62, 12).Value = 54156.668
$ws.Cells.Item(62, 13).Value = -1776
$ws.Cells.Item(62, 14).Value = -55404.668
$ws.Cells.Item(65, 8).Value = 41217.5
$ws.Cells.Item(65, 9).Value = 2400
$ws.Cells.Item(65, 10).Value = 54156.668
$ws.Cells.Item(65, 11).Value = 12000
$ws.Cells.Item(65, 12).Value = 270783.34
$ws.Cells.Item(65, 13).Value = -8880
$ws.Cells.Item(65, 14).Value = -277023.34
